# ---------------------------------------------------------------------------
# This script reproduces the commit:
#   "Commit on Wednesday. ... item_description is not being recognized and
#    it is breaking overrides"
#
# Summary of the change:
#  1) The "model" sheet gets a new "item_description" (string) row inserted
#     just before "acknowledge_auth", and the two rows that used to be typed
#     as "boolean" (is_distributed / is_override) are re-typed as "string".
#  2) A brand new "properties" worksheet is appended at the end of the
#     workbook, containing a small partition/aspect/key/type/value table.
#  3) The newly active sheet becomes "properties" (last tab).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Fix up the "model" sheet
# ---------------------------------------------------------------------------
$model = $wb.Worksheets.Item("model")
$model.Activate()

# Remove the two rows that reference the "boolean" / "is_distributed" /
# "is_override" shared strings entirely (not just clear them) so that those
# strings drop out of the shared-string table once unreferenced.
$model.Range("A10:B11").Delete() | Out-Null

# Push the old row 9 (acknowledge_auth) down to row 10, carrying its
# original formatting with it.
$model.Range("A9:B9").Copy($model.Range("A10")) | Out-Null

# Row 9 becomes the new item_description entry.
$model.Range("A9").Value = "string"
$model.Range("B9").Value = "item_description"

# Re-create is_distributed (row 11) / is_override (row 12) as plain strings,
# re-using the formatting already used by the other "string" rows above.
$model.Range("A8:B8").Copy($model.Range("A11")) | Out-Null
$model.Range("A8:B8").Copy($model.Range("A12")) | Out-Null
$model.Range("A11").Value = "string"
$model.Range("B11").Value = "is_distributed"
$model.Range("A12").Value = "string"
$model.Range("B12").Value = "is_override"

# Update the view: selection now sits on B11, tab no longer "selected"
# (that flag moves to the new last sheet once it is activated below).
$model.Range("B11").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2) Add the new "properties" worksheet as the last tab
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$properties = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$properties.Name = "properties"

# Helper cell used purely as a formatting donor so that the new font-8 based
# styles line up with the ones already present in the workbook instead of
# creating duplicate font/style entries.
$donor = $wb.Worksheets.Item("calculates").Range("B2")

function Set-PropCell($rng, [string]$text, [bool]$center) {
    $donor.Copy() | Out-Null
    $rng.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $rng.Value = $text
    $rng.Font.Color = 0
    if ($center) {
        $rng.HorizontalAlignment = -4108  # xlCenter
        $rng.VerticalAlignment = -4108    # xlCenter
    }
}

# Row 1 - headers
Set-PropCell $properties.Range("A1") "partition" $false
Set-PropCell $properties.Range("B1") "aspect" $false
Set-PropCell $properties.Range("C1") "key" $false
Set-PropCell $properties.Range("D1") "type" $false
Set-PropCell $properties.Range("E1") "value" $true

# Row 2
Set-PropCell $properties.Range("A2") "Table" $false
Set-PropCell $properties.Range("B2") "security" $false
Set-PropCell $properties.Range("C2") "unverifiedUserCanCreate" $false
Set-PropCell $properties.Range("D2") "boolean" $false
Set-PropCell $properties.Range("E2") "'false" $true

# Row 3
Set-PropCell $properties.Range("A3") "Table" $false
Set-PropCell $properties.Range("B3") "security" $false
Set-PropCell $properties.Range("C3") "filterTypeOnCreation" $false
Set-PropCell $properties.Range("D3") "string" $false
Set-PropCell $properties.Range("E3") "READ_ONLY" $true

# Row 4
Set-PropCell $properties.Range("A4") "Table" $false
Set-PropCell $properties.Range("B4") "security" $false
Set-PropCell $properties.Range("C4") "locked" $false
Set-PropCell $properties.Range("D4") "boolean" $false
Set-PropCell $properties.Range("E4") "'true" $true

# Row 5
Set-PropCell $properties.Range("A5") "FormType" $false
Set-PropCell $properties.Range("B5") "default" $false
Set-PropCell $properties.Range("C5") "FormType.formType" $false
Set-PropCell $properties.Range("D5") "string" $false
Set-PropCell $properties.Range("E5") "SURVEY" $false

# Row 6
Set-PropCell $properties.Range("A6") "SurveyUtil" $false
Set-PropCell $properties.Range("B6") "default" $false
Set-PropCell $properties.Range("C6") "SurveyUtil.formId" $false
Set-PropCell $properties.Range("D6") "string" $false
Set-PropCell $properties.Range("E6") "wrong_form" $false

# Select the whole table, matching the saved view state.
$properties.Range("A1:E6").Select() | Out-Null
$properties.Range("E6").Activate() | Out-Null

# Make "properties" the active tab (this is what's reflected as activeTab
# in the workbook view once saved).
$properties.Activate()
